$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.380.76"
$ws.Range("E2").Value = "  +7.60%  "
$ws.Range("D3").Value = "3.577.14"
$ws.Range("E3").Value = "  +3.36%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'417.68"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'130.52"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.770.12"
$ws.Range("E7").Value = "  +8.94%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.652"
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.785"
$ws.Range("E10").Value = "  +7.67%  "
$ws.Range("D11").Value = "'0.183"
$ws.Range("E11").Value = "  +29.78%  "
$ws.Range("D12").Value = "'0.0000341"
$ws.Range("E12").Value = "  +55.80%  "
$ws.Range("D13").Value = "'42.82"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").Value = "4.117.63"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'20.23"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "3.582.12"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("E19").Value = "  +4.76%  "
$ws.Range("D20").Value = "'12.52"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").Value = "67.206.58"
$ws.Range("E21").Value = "  +7.25%  "
$ws.Range("D22").Value = "'462.63"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "'90.25"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'3.18"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "'13.02"
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'35.28"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'9.91"
$ws.Range("E28").Value = "  -5.19%  "
$ws.Range("D29").Value = "'4.85"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'2.79"
$ws.Range("E30").Value = "  +5.54%  "
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("E32").Value = "  +3.95%  "
$ws.Range("D33").Value = "'7.32"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").Value = "'0.157"
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'39.31"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").Value = "'56.54"
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").Value = "0.0₃0791"
$ws.Range("E38").Value = "  +41.81%  "
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  +10.29%  "
$ws.Range("D41").Value = "'0.994"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'148.60"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.98"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.73"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("D47").Value = "'0.307"
$ws.Range("E47").Value = "  -5.03%  "
$ws.Range("D48").Value = "'1.97"
$ws.Range("E48").Value = "  -4.83%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'121.06"
$ws.Range("E49").Value = "  +10.26%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -4.80%  "
$ws.Range("E51").Value = "  +10.39%  "
